$p = $ppt.ActivePresentation
$p.Slides.Item(8).Delete()
